# This workbook ("Fruta, Terminal La Palmera de La Serena - Kiwi") receives a
# new weekly price record. The new record is inserted as row 451; every
# existing record from the old row 451 down to the old row 499 is pushed down
# by one row (old row 451 -> new row 452, ..., old row 499 -> new row 500).
#
# Net effect vs. the starting workbook:
#   - dimension grows from A1:T499 to A1:T500
#   - a new row 451 is created with the new weekly record
#   - every other data row (old 451..499) moves down by exactly one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 451; Excel shifts rows 451-499 down
# to 452-500 and extends the used range / dimension automatically.
$ws.Rows.Item(451).Insert()

# Populate the newly inserted row 451 with the new weekly record. Columns
# A, B, C, E, F, G, H, I, J, K and R carry the same constant values as every
# other row in this subset (market / product metadata).
$ws.Cells.Item(451, 1).Value = 8
$ws.Cells.Item(451, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(451, 3).Value = "Coquimbo"
$ws.Cells.Item(451, 4).Value = 45005
$ws.Cells.Item(451, 5).Value = 4
$ws.Cells.Item(451, 6).Value = "Fruta"
$ws.Cells.Item(451, 7).Value = 100101
$ws.Cells.Item(451, 8).Value = "Berries"
$ws.Cells.Item(451, 9).Value = 100101007
$ws.Cells.Item(451, 10).Value = "Kiwi"
$ws.Cells.Item(451, 11).Value = "Hayward"
$ws.Cells.Item(451, 12).Value = "Primera"
$ws.Cells.Item(451, 13).Value = 16
$ws.Cells.Item(451, 14).Value = 300000
$ws.Cells.Item(451, 15).Value = 310000
$ws.Cells.Item(451, 16).Value = 305000
$ws.Cells.Item(451, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(451, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(451, 19).Value = 678
$ws.Cells.Item(451, 20).Value = 450
